# edit.ps1 - Applies the row-data reshuffle described by the diff.
# Species observation records (columns A-AC) move between rows 2-17;
# shared columns (location/date/observer, P.. onward) are unaffected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Preserve "present-but-empty" Alder-Stadium/Kon/Metod cells (K, L, N)
#    by copying them before their source rows are overwritten below.
$ws.Range("K4").Copy($ws.Range("K3"))
$ws.Range("L4").Copy($ws.Range("L3"))
$ws.Range("N4").Copy($ws.Range("N3"))
$ws.Range("K10").Copy($ws.Range("K8"))
$ws.Range("L10").Copy($ws.Range("L8"))
$ws.Range("N10").Copy($ws.Range("N8"))

# 2) Write the relocated / updated values for each row.
$ws.Range("A2").Value = 104593632
$ws.Range("Q2").Value = 473815.7661137963
$ws.Range("R2").Value = 7013977.153526685
$ws.Range("A3").Value = 104593622
$ws.Range("B3").Value = 56395
$ws.Range("D3").Value = 'NT'
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = 'Tretåig hackspett'
$ws.Range("G3").Value = 'Picoides tridactylus'
$ws.Range("H3").Value = '(Linnaeus, 1758)'
$ws.Range("M3").Value = 'äldre spår'
$ws.Range("Q3").Value = 473830.905343441
$ws.Range("R3").Value = 7013897.29666794
$ws.Range("AC3").Value = 'ringhack'
$ws.Range("A4").Value = 104593631
$ws.Range("B4").Value = 96334
$ws.Range("D4").Value = 'VU'
$ws.Range("E4").Value = 220787
$ws.Range("F4").Value = 'Knärot'
$ws.Range("G4").Value = 'Goodyera repens'
$ws.Range("H4").Value = '(L.) R. Br.'
$ws.Range("Q4").Value = 473812.0075608135
$ws.Range("R4").Value = 7013958.714830574
$ws.Range("A5").Value = 104593630
$ws.Range("Q5").Value = 473798.8866438381
$ws.Range("R5").Value = 7013953.866335354
$ws.Range("A6").Value = 104593624
$ws.Range("Q6").Value = 473801.0947980214
$ws.Range("R6").Value = 7013892.583679659
$ws.Range("A7").Value = 104593620
$ws.Range("B7").Value = 89392
$ws.Range("D7").Value = 'NT'
$ws.Range("E7").Value = 1202
$ws.Range("F7").Value = 'Ullticka'
$ws.Range("G7").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H7").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q7").Value = 473808.0725733605
$ws.Range("R7").Value = 7013974.062789564
$ws.Range("A8").Value = 104593623
$ws.Range("B8").Value = 56395
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = 'Tretåig hackspett'
$ws.Range("G8").Value = 'Picoides tridactylus'
$ws.Range("H8").Value = '(Linnaeus, 1758)'
$ws.Range("M8").Value = 'äldre spår'
$ws.Range("Q8").Value = 473722.3678416939
$ws.Range("R8").Value = 7013918.902128431
$ws.Range("AC8").Value = 'ringhack'
$ws.Range("A9").Value = 104593627
$ws.Range("B9").Value = 96334
$ws.Range("D9").Value = 'VU'
$ws.Range("E9").Value = 220787
$ws.Range("F9").Value = 'Knärot'
$ws.Range("G9").Value = 'Goodyera repens'
$ws.Range("H9").Value = '(L.) R. Br.'
$ws.Range("Q9").Value = 473701.5160585373
$ws.Range("R9").Value = 7013906.458910029
$ws.Range("A10").Value = 104593628
$ws.Range("B10").Value = 96334
$ws.Range("D10").Value = 'VU'
$ws.Range("E10").Value = 220787
$ws.Range("F10").Value = 'Knärot'
$ws.Range("G10").Value = 'Goodyera repens'
$ws.Range("H10").Value = '(L.) R. Br.'
$ws.Range("Q10").Value = 473726.4678040863
$ws.Range("R10").Value = 7013923.823926651
$ws.Range("A11").Value = 104593626
$ws.Range("Q11").Value = 473718.6013391476
$ws.Range("R11").Value = 7013899.562304306
$ws.Range("A12").Value = 104593625
$ws.Range("Q12").Value = 473775.8828205758
$ws.Range("R12").Value = 7013898.645042086
$ws.Range("A13").Value = 104593629
$ws.Range("Q13").Value = 473760.983223469
$ws.Range("R13").Value = 7013952.372943264
$ws.Range("A14").Value = 104593637
$ws.Range("Q14").Value = 473780.3094452888
$ws.Range("R14").Value = 7013777.427834951
$ws.Range("A15").Value = 104593621
$ws.Range("B15").Value = 89392
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 1202
$ws.Range("F15").Value = 'Ullticka'
$ws.Range("G15").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H15").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Q15").Value = 473766.1571259646
$ws.Range("R15").Value = 7013701.408301079
$ws.Range("A16").Value = 104593634
$ws.Range("Q16").Value = 473769.3477768434
$ws.Range("R16").Value = 7013705.43688098
$ws.Range("A17").Value = 104593636
$ws.Range("Q17").Value = 473782.4035598941
$ws.Range("R17").Value = 7013757.588904253

# 3) Clear cells whose content moved away and has no replacement in this row.
$ws.Range("K4").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("AC4").ClearContents()
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()
$ws.Range("AC10").ClearContents()
